$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 11.44237473425257
$ws.Range("C2").Value = 3.704411347280387
$ws.Range("D2").Value = 9.552705627693321
$ws.Range("E2").Value = 11.69700795633191
$ws.Range("F2").Value = 50.68607110340694
$ws.Range("K2").Value = 10.97568907855668
$ws.Range("M2").Value = 14.84782912126718
$ws.Range("B3").Value = 11.34362332433227
$ws.Range("C3").Value = 3.649666658397357
$ws.Range("D3").Value = 9.388250751436674
$ws.Range("E3").Value = 11.4271104078
$ws.Range("F3").Value = 49.37148576754086
$ws.Range("K3").Value = 10.94998322847086
$ws.Range("M3").Value = 14.75478778018307
$ws.Range("B4").Value = 11.28818494904763
$ws.Range("C4").Value = 3.619505316429356
$ws.Range("D4").Value = 9.284724299055879
$ws.Range("E4").Value = 11.26136667513819
$ws.Range("F4").Value = 48.54807847348761
$ws.Range("K4").Value = 10.93930522970984
$ws.Range("M4").Value = 14.70243272264884
$ws.Range("B5").Value = 11.26692793595285
$ws.Range("C5").Value = 3.608107096826746
$ws.Range("D5").Value = 9.241918013633013
$ws.Range("E5").Value = 11.1939170698023
$ws.Range("F5").Value = 48.20882126946
$ws.Range("K5").Value = 10.93623994860149
$ws.Range("M5").Value = 14.68231567461216
$ws.Range("B6").Value = 11.26347958681193
$ws.Range("C6").Value = 3.60626899660653
$ws.Range("D6").Value = 9.234773397616019
$ws.Range("E6").Value = 11.18272562102473
$ws.Range("F6").Value = 48.15227499772809
$ws.Range("K6").Value = 10.93580866534142
$ws.Range("M6").Value = 14.6790493008559
$ws.Range("B7").Value = 11.28789283126425
$ws.Range("C7").Value = 3.61934795157548
$ws.Range("D7").Value = 9.284149471485133
$ws.Range("E7").Value = 11.26045652012322
$ws.Range("F7").Value = 48.54351765333974
$ws.Range("K7").Value = 10.93925868147933
$ws.Range("M7").Value = 14.7021564632805
$ws.Range("B8").Value = 11.40726485426869
$ws.Range("C8").Value = 3.684834907120255
$ws.Range("D8").Value = 9.496544426934099
$ws.Range("E8").Value = 11.60401101386664
$ws.Range("F8").Value = 50.23639124676016
$ws.Range("K8").Value = 10.96576801927881
$ws.Range("M8").Value = 14.81476934179019
$ws.Range("B9").Value = 11.68117845573765
$ws.Range("C9").Value = 3.839373611735051
$ws.Range("D9").Value = 9.891854797306292
$ws.Range("E9").Value = 12.27332883099231
$ws.Range("F9").Value = 53.41194598369502
$ws.Range("K9").Value = 11.05806812180469
$ws.Range("M9").Value = 15.0725769380567
$ws.Range("B10").Value = 11.9046271383813
$ws.Range("C10").Value = 3.96692395228656
$ws.Range("D10").Value = 10.16815792730414
$ws.Range("E10").Value = 12.75705989669339
$ws.Range("F10").Value = 55.63838383000779
$ws.Range("K10").Value = 11.15006230053062
$ws.Range("M10").Value = 15.28318677714536
$ws.Range("B11").Value = 12.01061195839896
$ws.Range("C11").Value = 4.027546042171851
$ws.Range("D11").Value = 10.2905510022568
$ws.Range("E11").Value = 12.97434257127292
$ws.Range("F11").Value = 56.62460640975666
$ws.Range("K11").Value = 11.19703857565526
$ws.Range("M11").Value = 15.38327779852169
$ws.Range("B12").Value = 12.05132598050291
$ws.Range("C12").Value = 4.050837442484678
$ws.Range("D12").Value = 10.3364052369354
$ws.Range("E12").Value = 13.05614419235969
$ws.Range("F12").Value = 56.99397989460226
$ws.Range("K12").Value = 11.21555132962569
$ws.Range("M12").Value = 15.42176612932296
$ws.Range("B13").Value = 12.04253240573798
$ws.Range("C13").Value = 4.045806931880029
$ws.Range("D13").Value = 10.32655194089188
$ws.Range("E13").Value = 13.03854925731611
$ws.Range("F13").Value = 56.91461444482972
$ws.Range("K13").Value = 11.21153231607284
$ws.Range("M13").Value = 15.41345140434989
$ws.Range("B14").Value = 12.01395015858823
$ws.Range("C14").Value = 4.029455744211914
$ws.Range("D14").Value = 10.29433345058974
$ws.Range("E14").Value = 12.9810824728142
$ws.Range("F14").Value = 56.65507833461289
$ws.Range("K14").Value = 11.19854720726178
$ws.Range("M14").Value = 15.38643268753898
$ws.Range("B15").Value = 11.9965169053719
$ws.Range("C15").Value = 4.019482667857266
$ws.Range("D15").Value = 10.2745338982725
$ws.Range("E15").Value = 12.94581778671373
$ws.Range("F15").Value = 56.49556523547277
$ws.Range("K15").Value = 11.19068729102549
$ws.Range("M15").Value = 15.3699583524695
$ws.Range("B16").Value = 11.89778433558923
$ws.Range("C16").Value = 3.963011156249233
$ws.Range("D16").Value = 10.16009122350831
$ws.Range("E16").Value = 12.74279724088585
$ws.Range("F16").Value = 55.57337365974681
$ws.Range("K16").Value = 11.14709435834292
$ws.Range("M16").Value = 15.27672925963607
$ws.Range("B17").Value = 11.83829512793001
$ws.Range("C17").Value = 3.929005946691698
$ws.Range("D17").Value = 10.08902560570922
$ws.Range("E17").Value = 12.61748251575429
$ws.Range("F17").Value = 55.00063251619699
$ws.Range("K17").Value = 11.12165577830785
$ws.Range("M17").Value = 15.22061246531633
$ws.Range("B18").Value = 11.80448942701469
$ws.Range("C18").Value = 3.909694000716925
$ws.Range("D18").Value = 10.04784112990733
$ws.Range("E18").Value = 12.5451486308217
$ws.Range("F18").Value = 54.66871791209221
$ws.Range("K18").Value = 11.10750774178572
$ws.Range("M18").Value = 15.18874063831639
$ws.Range("B19").Value = 11.79311528715967
$ws.Range("C19").Value = 3.903198894301439
$ws.Range("D19").Value = 10.03384423031111
$ws.Range("E19").Value = 12.52061616459458
$ws.Range("F19").Value = 54.55591824287033
$ws.Range("K19").Value = 11.10280089442724
$ws.Range("M19").Value = 15.17801984721596
$ws.Range("B20").Value = 11.84458564063019
$ws.Range("C20").Value = 3.932600579904408
$ws.Range("D20").Value = 10.09662282638468
$ws.Range("E20").Value = 12.63084959505192
$ws.Range("F20").Value = 55.06186137736457
$ws.Range("K20").Value = 11.12431380401284
$ws.Range("M20").Value = 15.22654448361333
$ws.Range("B21").Value = 12.02233007408403
$ws.Range("C21").Value = 4.034249696575884
$ws.Range("D21").Value = 10.30381033816641
$ws.Range("E21").Value = 12.99797546315046
$ws.Range("F21").Value = 56.73142329920917
$ws.Range("K21").Value = 11.20234171921845
$ws.Range("M21").Value = 15.39435307952887
$ws.Range("B22").Value = 12.14185529005171
$ws.Range("C22").Value = 4.102618881582043
$ws.Range("D22").Value = 10.43633559575721
$ws.Range("E22").Value = 13.23509037659055
$ws.Range("F22").Value = 57.7986471658409
$ws.Range("K22").Value = 11.25754922720895
$ws.Range("M22").Value = 15.50742684017475
$ws.Range("B23").Value = 12.07776980841755
$ws.Range("C23").Value = 4.065964575106728
$ws.Range("D23").Value = 10.36587423524127
$ws.Range("E23").Value = 13.10882111983352
$ws.Range("F23").Value = 57.23132035639065
$ws.Range("K23").Value = 11.22770346045754
$ws.Range("M23").Value = 15.4467762389704
$ws.Range("B24").Value = 11.84174046483327
$ws.Range("C24").Value = 3.930974701666595
$ws.Range("D24").Value = 10.09318914281384
$ws.Range("E24").Value = 12.6248072318238
$ws.Range("F24").Value = 55.03418800606917
$ws.Range("K24").Value = 11.12311062402084
$ws.Range("M24").Value = 15.22386139967276
$ws.Range("B25").Value = 11.60301839804319
$ws.Range("C25").Value = 3.794959186903964
$ws.Range("D25").Value = 9.787310752556808
$ws.Range("E25").Value = 12.09329909037857
$ws.Range("F25").Value = 52.57033774658073
$ws.Range("K25").Value = 11.0288156146142
$ws.Range("M25").Value = 14.99900575517166
